$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C (the old "Rent?" column), shifting D->C, E->D, F->E, G->F.
$ws.Columns("C").Delete()

# Leave the freshly-deleted (now shifted) column C selected, matching the
# post-delete state Excel leaves the sheet in.
$ws.Columns("C").Select()
